$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5045503333333333
$ws.Range("H2").Value = 1.513651
$ws.Range("I2").Value = 0.004193610593465892
$ws.Range("J2").Value = 0.004193610593465892
$ws.Range("M2").Value = 3.798569
$ws.Range("N2").Value = 11.395707
$ws.Range("O2").Value = 0.3335453021983209
$ws.Range("P2").Value = 0.3335453021983208
$ws.Range("Q2").Value = 1.916569255139667
$ws.Range("R2").Value = 17.249123296257
$ws.Range("S2").Value = 0.001398759112699661
$ws.Range("T2").Value = 0.00139875911269966
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5045503333333333
$ws.Range("H3").Value = 1.513651
$ws.Range("I3").Value = 0.004193610593465892
$ws.Range("J3").Value = 0.004193610593465892
$ws.Range("O3").Value = 0.3943739857244443
$ws.Range("P3").Value = 0.3943739857244442
$ws.Range("Q3").Value = 2.266094143988111
$ws.Range("R3").Value = 20.39484729589299
$ws.Range("S3").Value = 0.001653850924321396
$ws.Range("T3").Value = 0.001653850924321396
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5045503333333333
$ws.Range("H4").Value = 1.513651
$ws.Range("I4").Value = 0.004193610593465892
$ws.Range("J4").Value = 0.004193610593465892
$ws.Range("M4").Value = 3.098581666666667
$ws.Range("N4").Value = 9.295745
$ws.Range("O4").Value = 0.2720807120772349
$ws.Range("P4").Value = 0.2720807120772348
$ws.Range("Q4").Value = 1.563390412777222
$ws.Range("R4").Value = 14.070513714995
$ws.Range("S4").Value = 0.001141000556444836
$ws.Range("T4").Value = 0.001141000556444835
$ws.Range("I5").Value = 0.9882354706816875
$ws.Range("J5").Value = 0.9882354706816875
$ws.Range("M5").Value = 3.798569
$ws.Range("N5").Value = 11.395707
$ws.Range("O5").Value = 0.3335453021983209
$ws.Range("P5").Value = 0.3335453021983208
$ws.Range("Q5").Value = 451.6446336000043
$ws.Range("R5").Value = 4064.801702400039
$ws.Range("S5").Value = 0.3296212987116233
$ws.Range("T5").Value = 0.3296212987116233
$ws.Range("I6").Value = 0.9882354706816875
$ws.Range("J6").Value = 0.9882354706816875
$ws.Range("O6").Value = 0.3943739857244443
$ws.Range("P6").Value = 0.3943739857244442
$ws.Range("S6").Value = 0.3897343614070093
$ws.Range("T6").Value = 0.3897343614070092
$ws.Range("I7").Value = 0.9882354706816875
$ws.Range("J7").Value = 0.9882354706816875
$ws.Range("M7").Value = 3.098581666666667
$ws.Range("N7").Value = 9.295745
$ws.Range("O7").Value = 0.2720807120772349
$ws.Range("P7").Value = 0.2720807120772348
$ws.Range("Q7").Value = 368.4171016825961
$ws.Range("R7").Value = 3315.753915143365
$ws.Range("S7").Value = 0.2688798105630549
$ws.Range("T7").Value = 0.2688798105630549
$ws.Range("G8").Value = 0.9108879999999999
$ws.Range("H8").Value = 2.732664
$ws.Range("I8").Value = 0.007570918724846665
$ws.Range("J8").Value = 0.007570918724846665
$ws.Range("M8").Value = 3.798569
$ws.Range("N8").Value = 11.395707
$ws.Range("O8").Value = 0.3335453021983209
$ws.Range("P8").Value = 0.3335453021983208
$ws.Range("Q8").Value = 3.460070919272
$ws.Range("R8").Value = 31.140638273448
$ws.Range("S8").Value = 0.002525244373997907
$ws.Range("T8").Value = 0.002525244373997906
$ws.Range("G9").Value = 0.9108879999999999
$ws.Range("H9").Value = 2.732664
$ws.Range("I9").Value = 0.007570918724846665
$ws.Range("J9").Value = 0.007570918724846665
$ws.Range("O9").Value = 0.3943739857244443
$ws.Range("P9").Value = 0.3943739857244442
$ws.Range("Q9").Value = 4.091084330461332
$ws.Range("R9").Value = 36.81975897415199
$ws.Range("S9").Value = 0.002985773393113606
$ws.Range("T9").Value = 0.002985773393113606
$ws.Range("G10").Value = 0.9108879999999999
$ws.Range("H10").Value = 2.732664
$ws.Range("I10").Value = 0.007570918724846665
$ws.Range("J10").Value = 0.007570918724846665
$ws.Range("M10").Value = 3.098581666666667
$ws.Range("N10").Value = 9.295745
$ws.Range("O10").Value = 0.2720807120772349
$ws.Range("P10").Value = 0.2720807120772348
$ws.Range("Q10").Value = 2.822460857186667
$ws.Range("R10").Value = 25.40214771468
$ws.Range("S10").Value = 0.002059900957735152
$ws.Range("T10").Value = 0.002059900957735151
